# Auto-generated Excel COM-interop script
# Applies the per-cell value updates captured in the authoritative diff
# (Sheets/Maduin_Profits.xlsx) across the 8 worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value2 = 340.70834
$ws.Range("I15").Value2 = 340.70834
$ws.Range("K15").Value2 = 1022.12502
$ws.Range("M15").Value2 = -853.1250200000001
$ws.Range("H93").Value2 = 50000
$ws.Range("J93").Value2 = 50000
$ws.Range("L93").Value2 = 50000
$ws.Range("N93").Value2 = -54992
$ws.Range("H102").Value2 = 0
$ws.Range("J102").Value2 = 0
$ws.Range("L102").Value2 = 0
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value2 = 19096.23
$ws.Range("I113").Value2 = 21824.8
$ws.Range("J113").Value2 = 10001
$ws.Range("K113").Value2 = 21824.8
$ws.Range("L113").Value2 = 10001
$ws.Range("M113").Value2 = -18570.8
$ws.Range("N113").Value2 = -16509
$ws.Range("H129").Value2 = 2458
$ws.Range("I129").Value2 = 484.77777
$ws.Range("K129").Value2 = 1454.33331
$ws.Range("M129").Value2 = 3545.66669

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 4638.8887
$ws.Range("I45").Value2 = 1750
$ws.Range("J45").Value2 = 5000
$ws.Range("K45").Value2 = 1750
$ws.Range("L45").Value2 = 5000
$ws.Range("M45").Value2 = -1373
$ws.Range("N45").Value2 = -5754
$ws.Range("H92").Value2 = 44494.5
$ws.Range("I92").Value2 = 44990
$ws.Range("J92").Value2 = 43999
$ws.Range("K92").Value2 = 44990
$ws.Range("L92").Value2 = 43999
$ws.Range("M92").Value2 = -42494
$ws.Range("N92").Value2 = -48991
$ws.Range("H95").Value2 = 100208
$ws.Range("J95").Value2 = 100208
$ws.Range("L95").Value2 = 100208
$ws.Range("N95").Value2 = -105700
$ws.Range("H96").Value2 = 40000
$ws.Range("J96").Value2 = 40000
$ws.Range("L96").Value2 = 40000
$ws.Range("N96").Value2 = -45492
$ws.Range("H104").Value2 = 200000000
$ws.Range("J104").Value2 = 200000000
$ws.Range("L104").Value2 = 200000000
$ws.Range("N104").Value2 = -200006988
$ws.Range("H122").Value2 = 2974.0908
$ws.Range("I122").Value2 = 2376.8572
$ws.Range("K122").Value2 = 7130.571599999999
$ws.Range("M122").Value2 = -4680.571599999999
$ws.Range("H132").Value2 = 370
$ws.Range("I132").Value2 = 307
$ws.Range("K132").Value2 = 921
$ws.Range("M132").Value2 = 1609

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1000
$ws.Range("I107").Value2 = 1000
$ws.Range("K107").Value2 = 1000
$ws.Range("M107").Value2 = 920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 501.4
$ws.Range("J16").Value2 = 437
$ws.Range("L16").Value2 = 437
$ws.Range("N16").Value2 = -1011
$ws.Range("H23").Value2 = 10166
$ws.Range("I23").Value2 = 10166
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 10166
$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = -9926
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value2 = 10166
$ws.Range("I27").Value2 = 10166
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 10166
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = -9974
$ws.Range("N27").ClearContents()
$ws.Range("H58").Value2 = 1354.8422
$ws.Range("I58").Value2 = 982.86664
$ws.Range("K58").Value2 = 982.86664
$ws.Range("M58").Value2 = -779.86664
$ws.Range("H113").Value2 = 501.4
$ws.Range("J113").Value2 = 437
$ws.Range("L113").Value2 = 437
$ws.Range("N113").Value2 = -4777
$ws.Range("H122").Value2 = 492.4
$ws.Range("I122").Value2 = 492.4
$ws.Range("K122").Value2 = 1477.2
$ws.Range("M122").Value2 = 972.8000000000002
$ws.Range("H136").Value2 = 1354.8422
$ws.Range("I136").Value2 = 982.86664
$ws.Range("K136").Value2 = 2948.59992
$ws.Range("M136").Value2 = -398.5999199999997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 9332
$ws.Range("J80").Value2 = 11665.667
$ws.Range("L80").Value2 = 34997.001
$ws.Range("N80").Value2 = -36869.001
$ws.Range("H83").Value2 = 9332
$ws.Range("J83").Value2 = 11665.667
$ws.Range("L83").Value2 = 104991.003
$ws.Range("N83").Value2 = -114351.003
$ws.Range("H87").Value2 = 99.5
$ws.Range("I87").Value2 = 99.5
$ws.Range("K87").Value2 = 298.5
$ws.Range("M87").Value2 = 949.5
$ws.Range("H90").Value2 = 99.5
$ws.Range("I90").Value2 = 99.5
$ws.Range("K90").Value2 = 895.5
$ws.Range("M90").Value2 = 5344.5
$ws.Range("H131").Value2 = 1041.697
$ws.Range("J131").Value2 = 1055.7587
$ws.Range("L131").Value2 = 3167.2761
$ws.Range("N131").Value2 = -13247.2761
$ws.Range("J140").Value2 = 3000
$ws.Range("L140").Value2 = 9000
$ws.Range("N140").Value2 = -19360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value2 = 200000000
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 200000000
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 200000000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value2 = -200000224
$ws.Range("H8").Value2 = 200000000
$ws.Range("I8").Value2 = 0
$ws.Range("J8").Value2 = 200000000
$ws.Range("K8").Value2 = 0
$ws.Range("L8").Value2 = 200000000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value2 = -200000278
$ws.Range("H11").Value2 = 6076285
$ws.Range("J11").Value2 = 8499
$ws.Range("L11").Value2 = 8499
$ws.Range("N11").Value2 = -8777
$ws.Range("H14").Value2 = 866063.5600000001
$ws.Range("I14").Value2 = 2375125
$ws.Range("J14").Value2 = 3742.7144
$ws.Range("K14").Value2 = 2375125
$ws.Range("L14").Value2 = 3742.7144
$ws.Range("M14").Value2 = -2374957
$ws.Range("N14").Value2 = -4078.7144
$ws.Range("H92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("N92").ClearContents()
$ws.Range("H102").Value2 = 1486.5834
$ws.Range("I102").Value2 = 871.3333
$ws.Range("K102").Value2 = 871.3333
$ws.Range("M102").Value2 = 750.6667
$ws.Range("H113").Value2 = 1405.8889
$ws.Range("I113").Value2 = 1406.5
$ws.Range("K113").Value2 = 1406.5
$ws.Range("M113").Value2 = 763.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value2 = 500
$ws.Range("J3").Value2 = 500
$ws.Range("L3").Value2 = 500
$ws.Range("N3").Value2 = -724
$ws.Range("H15").Value2 = 500
$ws.Range("J15").Value2 = 500
$ws.Range("L15").Value2 = 500
$ws.Range("N15").Value2 = -840
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 0
$ws.Range("K23").Value2 = 0
$ws.Range("M23").ClearContents()
$ws.Range("H31").Value2 = 6000
$ws.Range("I31").Value2 = 5500
$ws.Range("J31").Value2 = 7000
$ws.Range("K31").Value2 = 5500
$ws.Range("L31").Value2 = 7000
$ws.Range("M31").Value2 = -5252
$ws.Range("N31").Value2 = -7496
$ws.Range("H55").Value2 = 374.06668
$ws.Range("I55").Value2 = 110.125
$ws.Range("J55").Value2 = 675.7143
$ws.Range("K55").Value2 = 110.125
$ws.Range("L55").Value2 = 675.7143
$ws.Range("M55").Value2 = 62.875
$ws.Range("N55").Value2 = -1021.7143
$ws.Range("H104").Value2 = 0
$ws.Range("J104").Value2 = 0
$ws.Range("L104").Value2 = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value2 = 1741.8572
$ws.Range("I122").Value2 = 1712.25
$ws.Range("K122").Value2 = 5136.75
$ws.Range("M122").Value2 = -2686.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 5504.769
$ws.Range("I14").Value2 = 5300
$ws.Range("K14").Value2 = 5300
$ws.Range("M14").Value2 = -5132
$ws.Range("H26").Value2 = 2500
$ws.Range("J26").Value2 = 2500
$ws.Range("L26").Value2 = 2500
$ws.Range("N26").Value2 = -3086
$ws.Range("H34").Value2 = 5000
$ws.Range("J34").Value2 = 5000
$ws.Range("L34").Value2 = 5000
$ws.Range("N34").Value2 = -5406
$ws.Range("H100").Value2 = 6971624.5
$ws.Range("I100").Value2 = 9957678
$ws.Range("K100").Value2 = 19915356
$ws.Range("M100").Value2 = -19914815
$ws.Range("H132").Value2 = 1889.0667
$ws.Range("I132").Value2 = 1963.2858
$ws.Range("K132").Value2 = 5889.857400000001
$ws.Range("M132").Value2 = -3359.857400000001
